$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (this also updates the _xlnm._FilterDatabase defined name
# reference automatically, since it points at the sheet name).
$ws.Name = "template_v6-03"

# Update NOTELOC subfield header labels to drop the leading underscore.
# Re-assigning these cell values rebuilds the shared string table, which also
# drops the now-unused old "_NOTELOC..." entries and appends the new ones.
$ws.Range("AB1").Value = "NOTELOCcoordinates"
$ws.Range("AD1").Value = "NOTELOClocuncm"
$ws.Range("AG1").Value = "NOTELOCsurveysite"
$ws.Range("AI1").Value = "NOTELOCdirections"

# Update the matching hyperlink subaddresses (drop leading underscore).
foreach ($h in $ws.Hyperlinks) {
    if ($h.SubAddress -eq "_noteloccoordinates") {
        $h.SubAddress = "noteloccoordinates"
    } elseif ($h.SubAddress -eq "_noteloclocuncm") {
        $h.SubAddress = "noteloclocuncm"
    } elseif ($h.SubAddress -eq "_notelocsurveysite") {
        $h.SubAddress = "notelocsurveysite"
    } elseif ($h.SubAddress -eq "_notelocdirections") {
        $h.SubAddress = "notelocdirections"
    }
}

# Remove the internal template placeholder values 1,2,3,4 from row 2.
$ws.Range("AB2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AG2").ClearContents()
$ws.Range("AI2").ClearContents()

# Update the current selection shown in the sheet view.
$ws.Range("AI1").Select()
